$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) The stray "_GoBack" bookmark currently sits in an empty paragraph
#    right after "NO ACEPTO EN SUS TERMINOS Y CONDICIONES EL PRESENTE
#    AVISO DE PRIVACIDAD". It needs to move to the end of the
#    "CULIACAN, SINALOA A {fecha}." paragraph further down, so delete
#    it from its current location first.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) "CULIACAN, SINALOA A {vigencia}." -> "CULIACAN, SINALOA A {fecha}."
#    The placeholder text is split across two runs ("{vigencia" / "}.")
#    with identical run formatting; the edit keeps it split into two
#    runs ("{fecha}" / ".") with that same formatting, so rebuild the
#    two runs via a literal OOXML fragment (a plain text replace would
#    normalize/merge same-format runs into a single run).
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("SINALOA A {vigencia}.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$target = $d.Range($rng.End - 11, $rng.End)

$xml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
       "<w:r><w:rPr><w:rFonts w:ascii='Century Gothic' w:hAnsi='Century Gothic' w:cs='Arial'/><w:b/><w:sz w:val='12'/><w:szCs w:val='12'/></w:rPr><w:t>{fecha}</w:t></w:r>" +
       "<w:r><w:rPr><w:rFonts w:ascii='Century Gothic' w:hAnsi='Century Gothic' w:cs='Arial'/><w:b/><w:sz w:val='12'/><w:szCs w:val='12'/></w:rPr><w:t>.</w:t></w:r>" +
       "<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/>" +
       "</w:p>"
$target.InsertXML($xml)
